# Polling data update: new poll results came in, so the oldest poll row
# (row 83, MidDate 43176 / Galaxy poll) is removed and all subsequent rows
# shift up by one. View/selection state is also updated to reflect where
# the author was last working in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Remove the obsolete poll entry (row 83); everything below shifts up one row.
$ws.Rows("83:83").Delete()

# Restore/update the window split + selection state to match where the
# author ended up after the edit (bottom pane starts at row 58, with the
# active cell in that pane at B82).
$aw = $excel.ActiveWindow
$aw.SplitRow = 57
$aw.SplitColumn = 0

[void]$ws.Range("D1").Select()
[void]$ws.Range("B82").Select()

# Best-effort: reflect the updated workbook window geometry recorded by
# the author's Excel session.
$aw.Left = 2895
$aw.Top = 2745
$aw.Width = 12570
$aw.Height = 16365
